$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price series gets a new latest-week record inserted at row 9;
# every existing record from row 9 down shifts one row lower (row 9 -> 10,
# ..., row 25 -> 26). Insert a blank row at 9 so rows 9-25 shift to 10-26,
# preserving all of their original data and formatting.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44544
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112044
$ws.Range("G9").Value = "Perejil"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 900
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 950
$ws.Range("N9").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 475
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = "Hortaliza"
